# Update cryptos list values (Price column D, Volume(1h) column E)
# Values are written with a leading apostrophe to force Excel to treat them
# as text (avoids numeric reinterpretation, e.g. "13.70" -> 13.7, or
# "36.566.64" -> parse error/garbage), then the cell style is reset to
# "Normal" so no stray number-format/style gets attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'36.566.64"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = "'  +0.54%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'1.961.99"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'  -0.03%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'244.66"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = "'  +0.74%  "
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'0.617"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'  +0.16%  "
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'59.07"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'  +1.52%  "
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "'  -0.04%  "
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "'  +2.87%  "
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'  -2.16%  "
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'  -0.09%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'22.24"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "'  +3.42%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'2.251.01"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "'  +1.12%  "
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "'  +0.84%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'13.70"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "'  +0.62%  "
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "'  +0.69%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'1.955.56"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'  +0.49%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'36.509.00"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'  +0.53%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'69.91"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "'  +0.57%  "
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "'  -0.14%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'228.91"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'  +0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'5.06"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "'  +0.50%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'2.45"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "'  +0.84%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'2.36"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'  +3.17%  "
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'  +8.89%  "
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "'  +0.20%  "
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'160.11"
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "'  -0.84%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'19.43"
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = "'  -0.10%  "
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.Value = "'  +1.82%  "
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = "'  +0.71%  "
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'4.72"
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = "'  +1.47%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'0.0619"
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "'  -1.19%  "
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = "'  +0.73%  "
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.Value = "'  -0.17%  "
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = "'  +6.06%  "
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.Value = "'  -4.92%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'3.35"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = "'  +10.92%  "
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = "'  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.0984"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "'  +0.79%  "
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'  +1.17%  "
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "'  +0.28%  "
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = "'  +0.87%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'16.00"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "'  +0.33%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'1.364.67"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'  +1.09%  "
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'  +0.78%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'87.69"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'  +0.10%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'7.13"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'  +0.56%  "
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'  +0.76%  "
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'2.142.19"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'  +1.15%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'43.68"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "'  -3.32%  "
$c.Style = "Normal"
